# Ravi Bishnoi.xlsx edit
#  1. Rename the sheet from "Sheet1" to "Ravi Bishnoi"
#  2. Insert a new column A ("matchNo") shifting the existing columns
#     (teamName..result) one slot to the right (B..M)
#  3. Populate the new column A with the header "matchNo" and the
#     value "21st" for the single data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet
$ws.Name = "Ravi Bishnoi"

# 2. Insert a blank column before column A; this shifts B:L -> C:M etc.
$ws.Columns.Item(1).Insert()

# 3. Fill in the new "matchNo" column
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "21st"
